$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a brand-new paragraph right after the "2022年6月3日星期五"
#    paragraph (currently paragraph 5), containing the 端午节 text that
#    used to live in the final paragraph. It must end up as two runs
#    (same rPr) - use InsertXML with a full package fragment so the
#    engine keeps the runs separate instead of merging them.
# ------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter() | Out-Null

$duanwuXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>中雨，今天是农历五月初五，中国传统端午节</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>：端午节，这一天我们要吃粽子，赛龙舟。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(6).Range.InsertXML($duanwuXml)

# ------------------------------------------------------------------
# 2. Insert a further new paragraph after it with "2022年6月7日星期二",
#    reproducing the same 3-run split ("2" / "022" / "年6月7日星期二")
#    used by the other date paragraphs in the document.
# ------------------------------------------------------------------
$d.Paragraphs.Item(6).Range.InsertParagraphAfter() | Out-Null

$dateXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t>022</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>年6月7日星期二</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(7).Range.InsertXML($dateXml)

# ------------------------------------------------------------------
# 3. The paragraph that used to hold the 端午节 text (now the last
#    paragraph in the document) gets new wording describing the first
#    day of the college entrance exam. Find/Replace across the whole
#    old two-run text collapses it back down to a single run, which
#    is what the target document looks like.
# ------------------------------------------------------------------
$last = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$last.Find.Execute(
    "中雨，今天是农历五月初五，中国传统端午节：端午节，这一天我们要吃粽子，赛龙舟。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "晴，今天是高考第一天，上午考语文，下午考数学。", 2) | Out-Null
